$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.075.85"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.828.63"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.90"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6238"
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07435"
$ws.Range("E8").Value = "  -1.94%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2920"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.21"
$ws.Range("E10").Value = "  +2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07699"
$ws.Range("D12").Value = "1.824.48"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.003"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6668"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.46"
$ws.Range("E15").Value = "  -0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009377"
$ws.Range("E16").Value = "  -6.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.953"
$ws.Range("E17").Value = "  -2.23%  "
$ws.Range("D18").Value = "29.079.96"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "2.083.47"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.58"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "222.94"
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.117"
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.76"
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1391"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.484"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.86"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.489"
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05766"
$ws.Range("E30").Value = "  +9.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.152"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.122"
$ws.Range("E32").Value = "  +2.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.210"
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.827"
$ws.Range("E34").Value = "  -1.31%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7375"
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.136"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.668"
$ws.Range("E37").Value = "  -0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01768"
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.481"
$ws.Range("E41").Value = "  +1.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8913"
$ws.Range("E42").Value = "  -1.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.03"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5086"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07573"
$ws.Range("E49").Value = "  +14.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4051"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.979"
$ws.Range("E51").Value = "  +0.48%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.762"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.222.59"
$ws.Range("E39").Value = "  -1.98%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.977.04"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.87"
$ws.Range("E46").Value = "  +2.18%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000124"
$ws.Range("E47").Value = "  -1.78%  "
